$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column = 想去人数)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 90
$ws1.Range("F4").Value = 260
$ws1.Range("F5").Value = 151
$ws1.Range("F6").Value = 251
$ws1.Range("F7").Value = 201
$ws1.Range("F8").Value = 1917
$ws1.Range("F10").Value = 4487
$ws1.Range("F12").Value = 314

# Sheet "全部类型" (F column = 想去人数)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 90
$ws4.Range("F6").Value = 260
$ws4.Range("F7").Value = 151
$ws4.Range("F8").Value = 251
$ws4.Range("F9").Value = 201
$ws4.Range("F12").Value = 1917
$ws4.Range("F14").Value = 4488
$ws4.Range("F16").Value = 314

$wb.Save()
